$p = $ppt.ActivePresentation

# The deck's slide master currently uses the "Integral" theme (ppt/theme/theme2.xml).
# The commit swaps the presentation's design back to the plain "Office Theme"
# color palette (dk2/lt2/accent1-6/hlink/folHlink). dk1/lt1 are already identical
# between the two themes, and the font scheme / format scheme are byte-identical
# between "Integral" and "Office Theme", so only the color scheme entries below
# need to change.
#
# ThemeColorScheme.Colors(n).RGB stores the value as a packed 0xBBGGRR (OLE
# color), i.e. the reverse byte order of the "RRGGBB" hex seen in the OOXML
# <a:srgbClr val="RRGGBB"/>. The hex literals below are written as 0xBBGGRR so
# that the saved XML ends up with the intended "RRGGBB" value.

$t = $p.SlideMaster.Theme.ThemeColorScheme

$t.Colors(1).RGB  = 0x000000   # dk1      -> 000000 (unchanged)
$t.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF (unchanged)
$t.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$t.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$t.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$t.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$t.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5 (unchanged)
$t.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$t.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$t.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$t.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$t.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
